$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 3 (the "MuSCs" sending-cluster row) entirely - cells shift up,
# and the now-unused "MuSCs" shared string is dropped from sharedStrings.xml.
$ws.Rows(3).Delete()

# Update the remaining row 2 values to the new TPM-derived figures.
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 1.565239
$ws.Range("N2").Value = 4.695717
$ws.Range("Q2").Value = 0.06605517278533334
$ws.Range("R2").Value = 0.594496555068
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
